$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 3.710289666666667
$ws.Range("N2").Value = 11.130869
$ws.Range("O2").Value = 0.2419880648107798
$ws.Range("P2").Value = 0.2419880648107798
$ws.Range("Q2").Value = 1.772083815328889
$ws.Range("R2").Value = 15.94875433796
$ws.Range("S2").Value = 0.2419880648107798
$ws.Range("T2").Value = 0.2419880648107798

# Row 3 updates
$ws.Range("O3").Value = 0.3723903391035988
$ws.Range("P3").Value = 0.3723903391035988
$ws.Range("S3").Value = 0.3723903391035988
$ws.Range("T3").Value = 0.3723903391035988

# Row 4 updates
$ws.Range("M4").Value = 2.323136666666667
$ws.Range("N4").Value = 6.96941
$ws.Range("O4").Value = 0.1515168347388597
$ws.Range("P4").Value = 0.1515168347388597
$ws.Range("Q4").Value = 1.109561047155556
$ws.Range("R4").Value = 9.986049424400001
$ws.Range("S4").Value = 0.1515168347388597
$ws.Range("T4").Value = 0.1515168347388597

# Row 5 updates
$ws.Range("M5").Value = 3.589418666666667
$ws.Range("N5").Value = 10.768256
$ws.Range("O5").Value = 0.2341047613467618
$ws.Range("P5").Value = 0.2341047613467618
$ws.Range("Q5").Value = 1.714354214115556
$ws.Range("R5").Value = 15.42918792704
$ws.Range("S5").Value = 0.2341047613467618
$ws.Range("T5").Value = 0.2341047613467618
